$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Signature / sponsor block rewrite.
#    "Responsable de la organización:" -> "Patrocinador:"
#    Drop the old "[Nombre del sitio de taxis]" / "Firma de los responsables:"
#    / "Juan Adolfo Bustillos Alatorre" / "[Nombre del sitio de taxis]" lines,
#    rename the trailing "[Nombre del editor]" line to "Taxico" (keeping the
#    _GoBack bookmark that already lives in that paragraph), and drop the
#    blank paragraph that used to follow it.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("Responsable de la organización:", $true, $false, $false, $false, $false, $true, 1, $false, "Patrocinador:", 2) | Out-Null

$idx = 0
$anchorIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Patrocinador:") {
        $anchorIdx = $idx
        break
    }
}

# The four paragraphs right after "Patrocinador:" are no longer needed.
$pFrom = $d.Paragraphs($anchorIdx + 1)
$pTo = $d.Paragraphs($anchorIdx + 4)
$d.Range($pFrom.Range.Start, $pTo.Range.End).Delete()

# Rename "[Nombre del editor]" -> "Taxico" in place (bookmark stays put).
$idx2 = 0
$editorIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx2 = $idx2 + 1
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "[Nombre del editor]") {
        $editorIdx = $idx2
        break
    }
}
$pEditor = $d.Paragraphs($editorIdx)
$pEditor.Range.Find.Execute("[Nombre del editor]", $false, $false, $false, $false, $false, $true, 1, $false, "Taxico", 2) | Out-Null

# Drop the blank paragraph that used to trail the signature block.
$pBlank = $d.Paragraphs($editorIdx + 1)
$pBlank.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Drop the stale "lastRenderedPageBreak" rendering hint that precedes
#    "El actor ingresa sus credenciales..." (the one inside the grid table
#    stays untouched).
# ---------------------------------------------------------------------------

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*El actor ingresa sus credenciales*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $rTouch = $d.Range($target.Range.Start, $target.Range.Start + 3)
    $rTouch.Text = "El" + [char]32 + [char]32
    $rExtra = $d.Range($target.Range.Start + 2, $target.Range.Start + 3)
    $rExtra.Delete()
}
